# [FEATURE] Add new filter types and Update domainmodel
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# --- Insert 4 blank rows before the old totals row (110) so it becomes row 114,
#     copying formatting the same way Excel does on a plain row-insert. ---
$ws.Rows("110:113").Insert()

# --- Fill in the previously-empty row 108 with the new "Filterarten" task ---
#     (copy row 107's formats first, same layout/style as the rest of this task group)
$ws.Range("A107:K107").Copy()
$ws.Range("A108:K108").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(108, 8).Clear()

$ws.Cells.Item(108, 1).Value = 18
$ws.Cells.Item(108, 2).Value = "Konzeptuelles Design"
$ws.Cells.Item(108, 3).Value = "Content Map"
$ws.Cells.Item(108, 4).Value = "[FEATURE]"
$ws.Cells.Item(108, 5).Value = "Filterarten bestimmen und konzipieren"
$ws.Cells.Item(108, 6).Value = 44379
$ws.Cells.Item(108, 7).Value = 44359
$ws.Cells.Item(108, 9).Formula = "=ROUNDUP(((SUM(K108-J108)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(108, 10).Value = 0.625
$ws.Cells.Item(108, 11).Value = 0.70833333333333337

# --- Close out the "week" block ending at row 103 with weekly-sum cells ---
$ws.Range("L94").Copy()
$ws.Range("L103").PasteSpecial(-4122)
$ws.Cells.Item(103, 12).Formula = "=SUM(H96:I103)"
$ws.Cells.Item(103, 13).NumberFormat = "0.00"
$ws.Cells.Item(103, 13).Formula = "=SUM(L103+19.5)"
$excel.CutCopyMode = 0

# --- New rows below the (now shifted) totals row 114: budget / credits info ---
#     D115/F115 keep the plain "right aligned" style used by the blank filler rows;
#     E115/G115 are new formulas, explicitly left-aligned (new style #33).
$ws.Range("D109").Copy()
$ws.Range("D115").PasteSpecial(-4122)
$ws.Range("F115").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(115, 4).Value = "Budget: "
$ws.Cells.Item(115, 6).Value = "Budget: "

$ws.Cells.Item(115, 5).Formula = "=135-E114"
$ws.Cells.Item(115, 5).HorizontalAlignment = -4131
$ws.Cells.Item(115, 7).Formula = "=315-G114"
$ws.Cells.Item(115, 7).HorizontalAlignment = -4131

$ws.Cells.Item(116, 2).Value = "Kredits erreicht:"
$ws.Cells.Item(116, 3).Formula = "=ROUNDUP(C114/30, 0)"

$ws.Cells.Item(117, 2).Value = "Kredits erfordert:"
$ws.Cells.Item(117, 3).Value = 15

# --- Update the second Prefix dropdown validation range to cover the new rows ---
$ws.Range("D41:D113").Validation.Delete()
$ws.Range("D41:D113").Validation.Add(3, 1, 1, "=$N$3:$N$6")
$ws.Range("D41:D113").Validation.InCellDropdown = $true
$ws.Range("D41:D113").Validation.ErrorTitle = "Prefix nicht unterstützt"
$ws.Range("D41:D113").Validation.ErrorMessage = "Es konnte kein korrekter Prefix ausgegeben werden`n"
$ws.Range("D41:D113").Validation.PromptTitle = "Prefix"
$ws.Range("D41:D113").Validation.InputMessage = "Wählen Sie einen Prefix aus"

# --- Restore selection near the edited area ---
$ws.Activate()
$ws.Range("M108").Select()
